# Add a new "2022-Q3" quarter sheet right after the "总计" (total) sheet,
# fill it with the new quarter's fund-holding detail data, and insert a new
# summary row at the top of the "总计" sheet's data table for 2022-Q3.
# All of the other existing quarter sheets just shift right by one tab
# position - their content is untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new summary row (2022-Q3) right after the
#    header row, so it appears as the newest quarter at the top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Use a template cell's format for the new row's index cell (column A),
# then overwrite every data row's values explicitly (top to bottom) with
# their final, literal target values - avoids float drift that a
# row-shift/Insert() would otherwise introduce on unrelated rows.
$total.Range("A2").Copy()
$total.Range("A9").PasteSpecial(-4122)

$totalData = @(
    @(0, "2022-Q3", 11, 0.16),
    @(1, "2022-Q2", 9, 0.46),
    @(2, "2022-Q1", 11, 2.01),
    @(3, "2021-Q4", 12, 2.1),
    @(4, "2021-Q3", 4, 0.18),
    @(5, "2021-Q2", 4, 1.02),
    @(6, "2021-Q1", 2, 0.1),
    @(7, "2020-Q4", 2, 0.12)
)

for ($i = 0; $i -lt $totalData.Length; $i++) {
    $r = $i + 2
    $row = $totalData[$i]
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------
# 2) Brand-new "2022-Q3" sheet, placed right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q3.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q3Data = @(
    @(0, "160135", "南方中证高铁产业指数（LOF）",   "1.84", "95.01", "2.61", "0.0480", 8),
    @(1, "006973", "太平睿盈混合A",                 "3.84", "28.79", "0.94", "0.0361", 3),
    @(2, "160639", "鹏华中证高铁产业指数（LOF）A",   "0.75", "94.62", "2.59", "0.0194", 8),
    @(3, "010157", "汇安中证500指数增强A",           "0.82", "91.29", "1.80", "0.0148", 6),
    @(4, "005599", "汇安量化优选灵活配置混合A",      "0.51", "94.40", "2.83", "0.0144", 9),
    @(5, "007669", "太平睿盈混合C",                  "1.04", "28.79", "0.94", "0.0098", 3),
    @(6, "010158", "汇安中证500指数增强C",           "0.53", "91.29", "1.80", "0.0095", 6),
    @(7, "007775", "汇安量化先锋混合A",              "0.23", "93.75", "3.10", "0.0071", 8),
    @(8, "007776", "汇安量化先锋混合C",              "0.11", "93.75", "3.10", "0.0034", 8),
    @(9, "015678", "鹏华中证高铁产业指数（LOF）C",   "0.06", "94.62", "2.59", "0.0016", 8),
    @(10, "005600", "汇安量化优选灵活配置混合C",     "0.02", "94.40", "2.83", "0.0006", 9)
)

# Fund codes (column B) must stay text (leading zeros matter), so force a
# text number format on that column before writing the values.
$q3.Range("B2:B12").NumberFormat = "@"

foreach ($row in $q3Data) {
    $r = [int]$row[0] + 2
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
}

# Match the look of the other quarter sheets: bold/centered/bordered
# header row + index column (style used by every other sheet's B1:H1/A-col).
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$total.Range("A2").Copy()
$q3.Range("A2:A12").PasteSpecial(-4122)

$q3.Range("A1").Select()
